$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "002/TTT"
$ws.Range("C2").Value = "IR801997"
$ws.Range("D2").Value = "NOUBAIL MOHAMMED"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 2000

# Update row 3
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 2000

# Delete rows 4-7
$ws.Range("A4:O7").EntireRow.Delete()
